# Rename AddressBook/Person-based diagram labels to their TaskManager/Task
# equivalents, as part of updating the other diagrams in the Developer Guide.

$p = $ppt.ActivePresentation

# --- Slide 3 -----------------------------------------------------------
$s3 = $p.Slides.Item(3)

# "deletePerson(p)" -> "deleteTask(p)"
$s3.Shapes.Item(16).TextFrame.TextRange.Characters(1, 12).Text = "deleteTask"

# "post(AddressBookChangedEvent)" -> "post(TaskManagerChangedEvent)"
# (Shape 17's box is narrow enough that the new, slightly-wider label makes
# PowerPoint's autofit wrap it onto a second line; restore the original
# single-line box height afterwards so only the label text changes.)
$tb32 = $s3.Shapes.Item(17)
$tb32Height = $tb32.Height
$tb32.TextFrame.TextRange.Characters(6, 23).Text = "TaskManagerChangedEvent"
$tb32.Height = $tb32Height + 0.00004

$s3.Shapes.Item(29).TextFrame.TextRange.Characters(6, 23).Text = "TaskManagerChangedEvent"

# "handleAddresssBookChangedEvent()" -> "handleTaskManagerChangedEvent()"
$s3.Shapes.Item(37).TextFrame.TextRange.Characters(1, 30).Text = "handleTaskManagerChangedEvent"
$s3.Shapes.Item(44).TextFrame.TextRange.Characters(1, 30).Text = "handleTaskManagerChangedEvent"

# --- Slide 4 -----------------------------------------------------------
$s4 = $p.Slides.Item(4)

# "PersonListPanel" -> "TaskListPanel"
$s4.Shapes.Item(11).TextFrame.TextRange.Characters(1, 15).Text = "TaskListPanel"

# "PersonCard" -> "TaskCard"
$s4.Shapes.Item(12).TextFrame.TextRange.Characters(1, 10).Text = "TaskCard"

# --- Slide 6 -----------------------------------------------------------
$s6 = $p.Slides.Item(6)

# "deletePerson(p)" -> "deleteTask(p)"
$s6.Shapes.Item(27).TextFrame.TextRange.Characters(1, 12).Text = "deleteTask"
